# Workbook was edited: two new weekly price rows for "Papa" (potato) at the
# Terminal Hortofrutícola Agro Chillán market were added to the dataset.
# The new records are inserted right before the existing row 325, which
# pushes all the following rows (old 325-414) down by two rows (new
# 327-416), exactly as shown in the unified diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 325-326 (this shifts old rows 325..414 down
# to 327..416, carrying along their formatting, same as in the diff).
$ws.Rows("325:326").Insert()

# Fixed columns shared by every record in this block.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$categoriaId = 100114001
$categoria = "Papa"
$unidadKg  = 25
$clasif    = "Hortaliza"

# --- New row 325: Asterix, 1a (guarda) ---
$ws.Cells.Item(325, 1).Value  = $mercadoId
$ws.Cells.Item(325, 2).Value  = $mercado
$ws.Cells.Item(325, 3).Value  = $region
$ws.Cells.Item(325, 4).Value  = 44841
$ws.Cells.Item(325, 5).Value  = $codreg
$ws.Cells.Item(325, 6).Value  = $categoriaId
$ws.Cells.Item(325, 7).Value  = $categoria
$ws.Cells.Item(325, 8).Value  = "Asterix"
$ws.Cells.Item(325, 9).Value  = "1a (guarda)"
$ws.Cells.Item(325, 10).Value = 120
$ws.Cells.Item(325, 11).Value = 6000
$ws.Cells.Item(325, 12).Value = 6500
$ws.Cells.Item(325, 13).Value = 6250
$ws.Cells.Item(325, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(325, 15).Value = "Región de Ñuble"
$ws.Cells.Item(325, 16).Value = 250
$ws.Cells.Item(325, 17).Value = $unidadKg
$ws.Cells.Item(325, 18).Value = $clasif

# --- New row 326: Patagonia, 1a (guarda) ---
$ws.Cells.Item(326, 1).Value  = $mercadoId
$ws.Cells.Item(326, 2).Value  = $mercado
$ws.Cells.Item(326, 3).Value  = $region
$ws.Cells.Item(326, 4).Value  = 44841
$ws.Cells.Item(326, 5).Value  = $codreg
$ws.Cells.Item(326, 6).Value  = $categoriaId
$ws.Cells.Item(326, 7).Value  = $categoria
$ws.Cells.Item(326, 8).Value  = "Patagonia"
$ws.Cells.Item(326, 9).Value  = "1a (guarda)"
$ws.Cells.Item(326, 10).Value = 120
$ws.Cells.Item(326, 11).Value = 6000
$ws.Cells.Item(326, 12).Value = 6500
$ws.Cells.Item(326, 13).Value = 6250
$ws.Cells.Item(326, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(326, 15).Value = "Región de Ñuble"
$ws.Cells.Item(326, 16).Value = 250
$ws.Cells.Item(326, 17).Value = $unidadKg
$ws.Cells.Item(326, 18).Value = $clasif
